$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Date: row 8, column B
$meta.Range("B8").Value = "2025-08-13T14:50:19+00:00"

# Count: row 22, column B (value is numeric-looking but must stay text,
# same as the existing "12" shared string -- force text storage, then
# restore the original (non-"@") cell formatting/style).
$meta.Range("B22").NumberFormat = "@"
$meta.Range("B22").Value = "13"
$meta.Range("B20").Copy()
$meta.Range("B22").PasteSpecial(-4122)

# --- Concepts sheet: append new concept row ---
$concepts = $wb.Worksheets.Item("Concepts")

# Clone row 13's formatting + values into row 14, then overwrite B14/C14.
$concepts.Range("A13:D13").Copy($concepts.Range("A14:D14"))

$concepts.Range("B14").Value = "Other"
$concepts.Range("C14").Value = "Other"
